$d = $word.ActiveDocument

# Disable "smart quotes" autocorrect so straight apostrophes in our
# replacement text aren't silently turned into curly ones.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# The document contains the same template repeated once per language
# (English, Portuguese, French, Thai, Vietnamese, Spanish). We must only
# touch the French copy, so first locate the "French" section heading
# paragraph and the paragraph right before the next language section
# begins, and scope every Find/Replace to that span.

$languageHeadings = @("English", "Portuguese", "French", "Thai", "Vietnamese", "Spanish")

$count = $d.Paragraphs.Count
$frenchStartPara = -1
$frenchEndPara = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($frenchStartPara -eq -1 -and $t -eq "French") {
        $frenchStartPara = $i
    } elseif ($frenchStartPara -ne -1 -and $frenchEndPara -eq -1 -and $languageHeadings -contains $t) {
        $frenchEndPara = $i - 1
    }
}
if ($frenchStartPara -eq -1) {
    throw "Could not locate the French section heading"
}
if ($frenchEndPara -eq -1) {
    $frenchEndPara = $count
}

$sectionStart = $d.Paragraphs.Item($frenchStartPara).Range.Start
$sectionEnd = $d.Paragraphs.Item($frenchEndPara).Range.End

function Replace-InSection([string]$findText, [string]$newText) {
    $rng = $d.Range($sectionStart, $sectionEnd)
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                                $true, 1, $false)
    if (-not $found) {
        throw "Could not find text: $findText"
    }
    # $rng now collapses to the matched text; set it directly so literal
    # characters (e.g. straight apostrophes) aren't auto-corrected by the
    # Find/Replace path.
    $rng.Text = $newText
}

# "Subject line" -> "Sujet" (bold run), and the run right after it gains
# a leading space before the colon:
# ": Rencontrez notre équipe a [CITY] | [DATE]" -> " : Rencontrez notre équipe a [CITY] | [DATE]"
Replace-InSection "Subject line" "Sujet"
Replace-InSection ": Rencontrez notre équipe a [CITY] | [DATE]" " : Rencontrez notre équipe a [CITY] | [DATE]"

# "Vous êtes invité à notre séminaire Deriv" -> "Vous êtes invité.e à notre séminaire Deriv"
Replace-InSection "Vous êtes invité à notre séminaire Deriv" "Vous êtes invité.e à notre séminaire Deriv"

# "Cher [PARTNER NAME], " -> "Cher [NOM DU PARTENAIRE], "
Replace-InSection "Cher [PARTNER NAME], " "Cher [NOM DU PARTENAIRE], "

# "Votre responsable national vous informera de l'emplacement exact d'ici [DATE]"
# -> "Votre responsable national vous informera du lieu exact d'ici [DATE]"
Replace-InSection "Votre responsable national vous informera de l'emplacement exact d'ici [DATE]" `
                   "Votre responsable national vous informera du lieu exact d'ici [DATE]"
